$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value() = "ID"
$ws.Range("B1").Value() = "NAME"
$ws.Range("C1").Value() = "ADDRESS"

Write-Host "done"
